$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "28.287.29"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value2 = "1.865.73"
$ws.Range("E3").Value = "  +3.04%  "

# Row 4
$ws.Range("D4").Value2 = "1.001"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'310.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("D6").Value2 = "1.001"
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
$ws.Range("D7").Value2 = "0.4983"
$ws.Range("E7").Value = "  -3.34%  "

# Row 8
$ws.Range("D8").Value2 = "0.3976"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.1000"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +27.74%  "

# Row 10
$ws.Range("E10").Value = "  +0.72%  "

# Row 11
$ws.Range("D11").Value2 = "41.35"
$ws.Range("E11").Value = "  +0.78%  "

# Row 12
$ws.Range("D12").Value2 = "6.467"
$ws.Range("E12").Value = "  +2.05%  "

# Row 13
$ws.Range("D13").Value2 = "20.87"
$ws.Range("E13").Value = "  +2.10%  "

# Row 14
$ws.Range("D14").Value2 = "1.859.29"
$ws.Range("E14").Value = "  +2.81%  "

# Row 15
$ws.Range("D15").Value2 = "1.001"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("D16").Value2 = "7.395"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17
$ws.Range("D17").Value = "'0.00001145"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.57%  "

# Row 18
$ws.Range("D18").Value2 = "93.59"
$ws.Range("E18").Value = "  +1.20%  "

# Row 19
$ws.Range("D19").Value2 = "0.06647"
$ws.Range("E19").Value = "  +1.20%  "

# Row 20
$ws.Range("D20").Value2 = "1.001"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21
$ws.Range("D21").Value2 = "17.37"
$ws.Range("E21").Value = "  +0.38%  "

# Row 22
$ws.Range("D22").Value2 = "6.077"
$ws.Range("E22").Value = "  +1.14%  "

# Row 23
$ws.Range("D23").Value2 = "28.401.13"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").Value2 = "11.34"
$ws.Range("E24").Value = "  +1.98%  "

# Row 25
$ws.Range("E25").Value = "  +1.00%  "

# Row 26
$ws.Range("D26").Value2 = "21.22"
$ws.Range("E26").Value = "  +3.40%  "

# Row 27
$ws.Range("D27").Value2 = "2.074.98"
$ws.Range("E27").Value = "  +2.72%  "

# Row 28
$ws.Range("D28").Value2 = "2.479"
$ws.Range("E28").Value = "  +2.34%  "

# Row 29
$ws.Range("D29").Value2 = "157.55"
$ws.Range("E29").Value = "  -2.03%  "

# Row 30
$ws.Range("D30").Value2 = "127.41"
$ws.Range("E30").Value = "  -0.25%  "

# Row 31
$ws.Range("D31").Value2 = "0.1057"
$ws.Range("E31").Value = "  -3.82%  "

# Row 32
$ws.Range("E32").Value = "  -1.21%  "

# Row 33
$ws.Range("D33").Value2 = "5.637"
$ws.Range("E33").Value = "  +1.06%  "

# Row 34
$ws.Range("D34").Value2 = "3.591"
$ws.Range("E34").Value = "  -1.94%  "

# Row 35
$ws.Range("D35").Value2 = "0.06816"
$ws.Range("E35").Value = "  -5.11%  "

# Row 36
$ws.Range("D36").Value2 = "9.239"
$ws.Range("E36").Value = "  +1.09%  "

# Row 37
$ws.Range("D37").Value2 = "0.02387"
$ws.Range("E37").Value = "  +1.13%  "

# Row 38
$ws.Range("D38").Value2 = "0.2166"
$ws.Range("E38").Value = "  -1.13%  "

# Row 39
$ws.Range("D39").Value2 = "5.025"
$ws.Range("E39").Value = "  -0.54%  "

# Row 40
$ws.Range("D40").Value = "'11.50"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.78%  "

# Row 41
$ws.Range("D41").Value2 = "0.6291"
$ws.Range("E41").Value = "  +1.64%  "

# Row 42
$ws.Range("D42").Value = "'1.180"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.93%  "

# Row 43
$ws.Range("D43").Value2 = "1.001"
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
$ws.Range("D44").Value2 = "13.38"
$ws.Range("E44").Value = "  +0.83%  "

# Row 45
$ws.Range("D45").Value2 = "0.5989"
$ws.Range("E45").Value = "  -0.11%  "

# Row 46
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value2 = "1.281"
$ws.Range("E46").Value = "  -1.50%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value2 = "3.677"
$ws.Range("E47").Value = "  -1.69%  "

# Row 48
$ws.Range("D48").Value2 = "124.94"
$ws.Range("E48").Value = "  -0.40%  "

# Row 49
$ws.Range("D49").Value2 = "1.982"
$ws.Range("E49").Value = "  +2.98%  "

# Row 50
$ws.Range("D50").Value2 = "1.191"
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.120"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.89%  "
